$d = $word.ActiveDocument

# Locate the paragraph currently holding "Problème d’optimisation ..." —
# it is identified by its distinctive text and carries the "_GoBack"
# bookmark in the original document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Problème d’optimisation*") {
        $target = $cand
        break
    }
}

# --- Turn this paragraph into the new "Requêtes SQL :" bullet ------------
$target.Range.Text = "Requêtes SQL :"

# --- Grow the new ilvl=1 sub-bullets after it -----------------------------
$target.Range.InsertParagraphAfter()
$p = $target.Next()
$p.Range.Text = "Cartes "
$p.Range.ListFormat.ListLevelNumber = 2

$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Blocs + éditions"
$p.Range.ListFormat.ListLevelNumber = 2

$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Types"
$p.Range.ListFormat.ListLevelNumber = 2

$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Sous-types"
$p.Range.ListFormat.ListLevelNumber = 2

$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Super-types (éventuellement à supprimer plus tard si non-nécessaire)"
$p.Range.ListFormat.ListLevelNumber = 2

$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Capacités"
$p.Range.ListFormat.ListLevelNumber = 2

# This paragraph holds the "Ids : ... symbol ... ids des cartes" bullet —
# type the whole text (with a placeholder symbol char) in a single
# assignment so every run starts life with the paragraph's sz=24, then
# re-style just the symbol character afterwards.
$p.Range.InsertParagraphAfter()
$idsPara = $p.Next()
$idsPara.Range.ListFormat.ListLevelNumber = 2
$prefixText = "Ids : Pour avoir la correspondance multiverseid "
$idsPara.Range.Text = $prefixText + "" + " ids des cartes"

# --- Re-insert the original "Problème d’optimisation" bullet -------------
$idsPara.Range.InsertParagraphAfter()
$p = $idsPara.Next()
$p.Range.Text = "Problème d’optimisation : 2 algos ont été implémentés indépendamment (1 par membre du binôme), ils seront présentés lors de la réunion"
$p.Range.ListFormat.ListLevelNumber = 1

# --- New closing bullet "Ce début de rapport =)" --------------------------
$p.Range.InsertParagraphAfter()
$p = $p.Next()
$p.Range.Text = "Ce début de rapport =)"
$p.Range.ListFormat.ListLevelNumber = 1

# --- Now restyle the placeholder symbol character as Wingdings -----------
$symPos = $idsPara.Range.Start + $prefixText.Length
$symRange = $d.Range($symPos, $symPos + 1)
$symRange.Font.Name = "Wingdings"

# --- Move the _GoBack bookmark onto the closing bullet --------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$endPos = $p.Range.End - 1
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Drop one of the two now-redundant empty trailing paragraphs ----------
$after = $p.Next()
$after.Range.Delete()

Write-Host "done"
